$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[58.87385302329284, 68.02345237381314]"
$ws.Range("T2").Value = "[46.51986363576425, 52.799700046534184]"
$ws.Range("L3").Value = "[56.23779105879879, 69.99615276776865]"
$ws.Range("T3").Value = "[45.54280939407556, 52.68035776707866]"
